# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right after "总计" (i.e. before "2022-Q2"),
#    populate it with the quarterly fund-holding data, matching the layout
#    used by the existing quarter sheets.
# 2. Insert a new row at the top of the "总计" (totals) sheet's data for the
#    2022-Q4 summary figures, pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook

function Set-TextCell($sheet, $row, $col, $val) {
    # Force the cell to be stored as text (inline/shared string) even when
    # the value looks numeric (fund codes with leading zeros, ratios like
    # "0.80" that must keep their trailing zero, etc.) - mirrors how the
    # source sheets store these columns.
    $c = $sheet.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value2 = $val
    $c.Style = "Normal"
}

function Set-NumCell($sheet, $row, $col, $val) {
    $sheet.Cells.Item($row, $col).Value2 = $val
}

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q4" sheet right after "总计"
# ---------------------------------------------------------------------

$totalSheet = $wb.Worksheets.Item(1)
$q4Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q4Sheet.Name = "2022-Q4"

# Match page setup / outline settings used by the rest of the workbook.
$q4Sheet.PageSetup.LeftMargin = 54
$q4Sheet.PageSetup.RightMargin = 54
$q4Sheet.PageSetup.TopMargin = 72
$q4Sheet.PageSetup.BottomMargin = 72
$q4Sheet.PageSetup.HeaderMargin = 36
$q4Sheet.PageSetup.FooterMargin = 36
$q4Sheet.Outline.SummaryRow = 1
$q4Sheet.Outline.SummaryColumn = 1

# Copy the header-row formatting (bold/border/center) from an existing
# quarter sheet so the new sheet matches the others exactly.
$formatSource = $wb.Worksheets.Item("2022-Q2")
$formatSource.Range("B1:H1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)

Set-TextCell $q4Sheet 1 2 "基金代码"
Set-TextCell $q4Sheet 1 3 "基金名称"
Set-TextCell $q4Sheet 1 4 "基金规模"
Set-TextCell $q4Sheet 1 5 "股票总仓位"
Set-TextCell $q4Sheet 1 6 "仓位占比"
Set-TextCell $q4Sheet 1 7 "持有市值(亿元)"
Set-TextCell $q4Sheet 1 8 "仓位排名"
# re-apply the header style/border that Set-TextCell resets via NumberFormat
$formatSource.Range("B1:H1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)

# Copy the "A column" (row index) style used on the other quarter sheets.
$formatSource.Cells.Item(2, 1).Copy()
$q4Sheet.Cells.Item(2, 1).PasteSpecial(-4122)
$q4Sheet.Cells.Item(3, 1).Copy()
$q4Sheet.Range("A3:A8").PasteSpecial(-4122)

Set-NumCell  $q4Sheet 2 1 0
Set-TextCell $q4Sheet 2 2 "004634"
Set-TextCell $q4Sheet 2 3 "新疆前海联合泳涛灵活配置混合A"
Set-TextCell $q4Sheet 2 4 "1.37"
Set-TextCell $q4Sheet 2 5 "92.47"
Set-TextCell $q4Sheet 2 6 "5.12"
Set-TextCell $q4Sheet 2 7 "0.0701"
Set-NumCell  $q4Sheet 2 8 9

Set-NumCell  $q4Sheet 3 1 1
Set-TextCell $q4Sheet 3 2 "006235"
Set-TextCell $q4Sheet 3 3 "东方城镇消费主题混合"
Set-TextCell $q4Sheet 3 4 "0.64"
Set-TextCell $q4Sheet 3 5 "90.44"
Set-TextCell $q4Sheet 3 6 "4.49"
Set-TextCell $q4Sheet 3 7 "0.0287"
Set-NumCell  $q4Sheet 3 8 9

Set-NumCell  $q4Sheet 4 1 2
Set-TextCell $q4Sheet 4 2 "007041"
Set-TextCell $q4Sheet 4 3 "新疆前海联合泳涛灵活配置混合C"
Set-TextCell $q4Sheet 4 4 "0.42"
Set-TextCell $q4Sheet 4 5 "92.47"
Set-TextCell $q4Sheet 4 6 "5.12"
Set-TextCell $q4Sheet 4 7 "0.0215"
Set-NumCell  $q4Sheet 4 8 9

Set-NumCell  $q4Sheet 5 1 3
Set-TextCell $q4Sheet 5 2 "002872"
Set-TextCell $q4Sheet 5 3 "华夏智胜价值成长股票C"
Set-TextCell $q4Sheet 5 4 "2.68"
Set-TextCell $q4Sheet 5 5 "92.97"
Set-TextCell $q4Sheet 5 6 "0.80"
Set-TextCell $q4Sheet 5 7 "0.0214"
Set-NumCell  $q4Sheet 5 8 6

Set-NumCell  $q4Sheet 6 1 4
Set-TextCell $q4Sheet 6 2 "009619"
Set-TextCell $q4Sheet 6 3 "博时女性消费主题混合A"
Set-TextCell $q4Sheet 6 4 "0.56"
Set-TextCell $q4Sheet 6 5 "72.72"
Set-TextCell $q4Sheet 6 6 "3.27"
Set-TextCell $q4Sheet 6 7 "0.0183"
Set-NumCell  $q4Sheet 6 8 8

Set-NumCell  $q4Sheet 7 1 5
Set-TextCell $q4Sheet 7 2 "002871"
Set-TextCell $q4Sheet 7 3 "华夏智胜价值成长股票A"
Set-TextCell $q4Sheet 7 4 "0.92"
Set-TextCell $q4Sheet 7 5 "92.97"
Set-TextCell $q4Sheet 7 6 "0.80"
Set-TextCell $q4Sheet 7 7 "0.0074"
Set-NumCell  $q4Sheet 7 8 6

Set-NumCell  $q4Sheet 8 1 6
Set-TextCell $q4Sheet 8 2 "009620"
Set-TextCell $q4Sheet 8 3 "博时女性消费主题混合C"
Set-TextCell $q4Sheet 8 4 "0.03"
Set-TextCell $q4Sheet 8 5 "72.72"
Set-TextCell $q4Sheet 8 6 "3.27"
Set-TextCell $q4Sheet 8 7 "0.0010"
Set-NumCell  $q4Sheet 8 8 8

# ---------------------------------------------------------------------
# Step 2: insert the 2022-Q4 summary row into the "总计" sheet
# ---------------------------------------------------------------------

$totalSheet.Rows.Item(2).Insert()

# restore the "row index" style (bold/border/center) on the new A2 cell
$totalSheet.Cells.Item(3, 1).Copy()
$totalSheet.Cells.Item(2, 1).PasteSpecial(-4122)

Set-NumCell $totalSheet 2 1 0
$totalSheet.Cells.Item(2, 2).Style = "Normal"
Set-NumCell $totalSheet 2 2 "2022-Q4"
$totalSheet.Cells.Item(2, 3).Style = "Normal"
Set-NumCell $totalSheet 2 3 7
$totalSheet.Cells.Item(2, 4).Style = "Normal"
Set-NumCell $totalSheet 2 4 0.17

# ---------------------------------------------------------------------
# Restore the originally-active sheet/selection.
# ---------------------------------------------------------------------
$totalSheet.Activate()
